# Scheduled-runner price/profit refresh for the Tiamat Profits workbook.
# Re-pushes freshly recomputed currentAveragePrice* / LevePrice* / LeveProfit*
# figures (columns H-N) for the leves whose market data moved since the last
# run. Values are literal (no formulas live in these sheets), so each refreshed
# cell is written directly; a couple of rows gain/lose a NQ or HQ profit cell
# entirely when that side of the recipe stops (or starts) applying.

$wb = $excel.ActiveWorkbook

# ===== ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 86: Filling in the Blanks
$ws.Range("H86").Value = 1707.6428
$ws.Range("I86").Value = 1100.5555
$ws.Range("J86").Value = 2800.4
$ws.Range("K86").Value = 1100.5555
$ws.Range("L86").Value = 2800.4
$ws.Range("M86").Value = 22.44450000000006
$ws.Range("N86").Value = -5046.4
# Row 89: Ink into Antiquity (L)
$ws.Range("H89").Value = 1707.6428
$ws.Range("I89").Value = 1100.5555
$ws.Range("J89").Value = 2800.4
$ws.Range("K89").Value = 5502.7775
$ws.Range("L89").Value = 14002
$ws.Range("M89").Value = 113.2224999999999
$ws.Range("N89").Value = -25234
# Row 129: Practical Command
$ws.Range("H129").Value = 8580.799999999999
$ws.Range("J129").Value = 10601
$ws.Range("L129").Value = 31803
$ws.Range("N129").Value = -41803
# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 22728748
$ws.Range("I132").Value = 1381.4166
$ws.Range("J132").Value = 125001896
$ws.Range("K132").Value = 4144.2498
$ws.Range("L132").Value = 375005688
$ws.Range("M132").Value = -1614.2498
$ws.Range("N132").Value = -375010748
# Row 138: All-night Crafting
$ws.Range("H138").Value = 1769.4684
$ws.Range("I138").Value = 1215.3636
$ws.Range("J138").Value = 1983.3334
$ws.Range("K138").Value = 3646.0908
$ws.Range("L138").Value = 5950.0002
$ws.Range("M138").Value = 1493.9092
$ws.Range("N138").Value = -16230.0002

# ===== ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 113222.2
$ws.Range("I132").Value = 169537.33
$ws.Range("K132").Value = 508611.99
$ws.Range("M132").Value = -506081.99

# ===== CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 6: Got Your Back
$ws.Range("H6").Value = 745072.3
$ws.Range("I6").Value = 5000000
$ws.Range("J6").Value = 35917.668
$ws.Range("K6").Value = 5000000
$ws.Range("L6").Value = 35917.668
$ws.Range("M6").Value = -4999887
$ws.Range("N6").Value = -36143.668
# Row 7: Gridania's Got Talent
$ws.Range("H7").Value = 59.875
$ws.Range("I7").Value = 84
$ws.Range("J7").Value = 48.909092
$ws.Range("K7").Value = 84
$ws.Range("L7").Value = 48.909092
$ws.Range("M7").Value = 29
$ws.Range("N7").Value = -274.909092
# Row 17: Say It with Spears
$ws.Range("H17").Value = 2632.5
$ws.Range("I17").Value = 2275
$ws.Range("J17").Value = 2990
$ws.Range("K17").Value = 2275
$ws.Range("L17").Value = 2990
$ws.Range("M17").Value = -2101
$ws.Range("N17").Value = -3338
# Row 25: Bowing to Necessity
$ws.Range("H25").Value = 35479.43
$ws.Range("I25").Value = 11000
$ws.Range("J25").Value = 39559.332
$ws.Range("K25").Value = 11000
$ws.Range("L25").Value = 39559.332
$ws.Range("M25").Value = -10826
$ws.Range("N25").Value = -39907.332
# Row 41: The Lone Bowman
$ws.Range("H41").Value = 8677.8125
$ws.Range("I41").Value = 4546.6665
$ws.Range("J41").Value = 11156.5
$ws.Range("K41").Value = 4546.6665
$ws.Range("L41").Value = 11156.5
$ws.Range("M41").Value = -4118.6665
$ws.Range("N41").Value = -12012.5
# Row 50: The Arsenal of Theocracy
$ws.Range("H50").Value = 16000
$ws.Range("J50").Value = 16000
$ws.Range("L50").Value = 16000
$ws.Range("N50").Value = -17250
# Row 51: Greenstone for Greenhorns
$ws.Range("H51").Value = 10000
$ws.Range("J51").Value = 10000
$ws.Range("L51").Value = 10000
$ws.Range("N51").Value = -11472
# Row 59: Bow Down to Magic
$ws.Range("H59").Value = 12097.143
$ws.Range("J59").Value = 12097.143
$ws.Range("L59").Value = 12097.143
$ws.Range("N59").Value = -14387.143
# Row 61: Incant Now, Think Later
$ws.Range("H61").Value = 10000
$ws.Range("J61").Value = 10000
$ws.Range("L61").Value = 10000
$ws.Range("N61").Value = -10696
# Row 68: Do You Even String Bow
$ws.Range("H68").Value = 23457
$ws.Range("I68").Value = 9200
$ws.Range("J68").Value = 32961.668
$ws.Range("K68").Value = 9200
$ws.Range("L68").Value = 32961.668
$ws.Range("M68").Value = -8451
$ws.Range("N68").Value = -34459.668
# Row 71: Win One Bow, Get Three Free (L)
$ws.Range("H71").Value = 23457
$ws.Range("I71").Value = 9200
$ws.Range("J71").Value = 32961.668
$ws.Range("K71").Value = 27600
$ws.Range("L71").Value = 98885.00399999999
$ws.Range("M71").Value = -23856
$ws.Range("N71").Value = -106373.004
# Row 74: License to Heal
$ws.Range("H74").Value = 12363
$ws.Range("J74").Value = 12363
$ws.Range("L74").Value = 12363
$ws.Range("N74").Value = -14111
# Row 77: Purified Polyrhythm (L)
$ws.Range("H77").Value = 12363
$ws.Range("J77").Value = 12363
$ws.Range("L77").Value = 37089
$ws.Range("N77").Value = -45825
# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 2255.138
$ws.Range("I132").Value = 1462.7894
$ws.Range("J132").Value = 3760.6
$ws.Range("K132").Value = 4388.3682
$ws.Range("L132").Value = 11281.8
$ws.Range("M132").Value = -1858.3682
$ws.Range("N132").Value = -16341.8

# ===== CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap
$ws.Range("H5").Value = 20801
$ws.Range("J5").Value = 20801
$ws.Range("L5").Value = 62403
$ws.Range("N5").Value = -62627
# Row 129: Comfort Food
$ws.Range("H129").Value = 39827.73
$ws.Range("I129").Value = 686.6667
$ws.Range("J129").Value = 51570.05
$ws.Range("K129").Value = 2060.0001
$ws.Range("L129").Value = 154710.15
$ws.Range("M129").Value = 2939.9999
$ws.Range("N129").Value = -164710.15
# Row 134: Don't Knock It Till You've Tried It
$ws.Range("H134").Value = 4370.16
$ws.Range("I134").Value = 1850.1333
$ws.Range("J134").Value = 8150.2
$ws.Range("K134").Value = 5550.3999
$ws.Range("L134").Value = 24450.6
$ws.Range("M134").Value = -480.3999000000003
$ws.Range("N134").Value = -34590.6
# Row 135: Not-so-secret Ingredient
$ws.Range("H135").Value = 20801
$ws.Range("J135").Value = 20801
$ws.Range("L135").Value = 187209
$ws.Range("N135").Value = -192279

# ===== GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 132: On Board for Lar
$ws.Range("H132").Value = 28543.395
$ws.Range("I132").Value = 1525.5834
$ws.Range("K132").Value = 4576.7502
$ws.Range("M132").Value = -2046.7502

# ===== LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 18: Simply the Best
$ws.Range("H18").Value = 16666.666
$ws.Range("I18").Value = 20000
$ws.Range("J18").Value = 10000
$ws.Range("K18").Value = 20000
$ws.Range("L18").Value = 10000
$ws.Range("M18").Value = -19828
$ws.Range("N18").Value = -10344
# Row 61: Spelling Me Softly
$ws.Range("H61").Value = 4772.273
$ws.Range("I61").Value = 3315
$ws.Range("K61").Value = 3315
$ws.Range("M61").Value = -3113
# Row 113: Peace in Rest
$ws.Range("H113").Value = 4772.273
$ws.Range("I113").Value = 3315
$ws.Range("K113").Value = 3315
$ws.Range("M113").Value = -1145
# Row 122: Hell on Leather
$ws.Range("H122").Value = 1617.375
$ws.Range("I122").Value = 1639.1428
$ws.Range("J122").Value = 1600.4445
$ws.Range("K122").Value = 4917.428400000001
$ws.Range("L122").Value = 4801.333500000001
$ws.Range("M122").Value = -2467.428400000001
$ws.Range("N122").Value = -9701.333500000001

# ===== WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 49: A Leg Up on the Cold
$ws.Range("H49").Value = 67374.664
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 67374.664
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 67374.664
$ws.Range("M49").ClearContents()  # was -9770, no longer applicable
$ws.Range("N49").Value = -67834.664
# Row 126: A Polished Purchase
$ws.Range("H126").Value = 1459.375
$ws.Range("I126").Value = 1210
$ws.Range("J126").Value = 1572.7273
$ws.Range("K126").Value = 3630
$ws.Range("L126").Value = 4718.1819
$ws.Range("M126").Value = -1160
$ws.Range("N126").Value = -9658.1819
# Row 132: Comfy Cabins
$ws.Range("H132").Value = 1457.7709
$ws.Range("I132").Value = 1180.5143
$ws.Range("K132").Value = 3541.5429
$ws.Range("M132").Value = -1011.5429
